# Apply the "cryptos list" refresh described by the commit.
# Every data row (2-51) gets new Price (D) / Volume(1h) (E) text,
# and rows 39/40 swap which coin (TheSandbox / InternetComputer) sits where.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.945.06"
$ws.Range("E2").Value = "  +1.26%  "

$ws.Range("D3").Value = "1.768.44"
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "

$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4562"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.47%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07384"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.93%  "

$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.008"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.188"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.64%  "

$ws.Range("D16").Value = "1.768.91"
$ws.Range("E16").Value = "  +0.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.94%  "

$ws.Range("E18").Value = "  +0.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06446"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.777"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.51%  "

$ws.Range("D23").Value = "27.969.73"
$ws.Range("E23").Value = "  +1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.102"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("D28").Value = "1.971.99"
$ws.Range("E28").Value = "  +0.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.175"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.079"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09244"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.609"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.659"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02282"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06129"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2093"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6265"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.921"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.182"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.381"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.820"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.735"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5855"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06824"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.35%  "
